# adding site info to manu
# Reorder the sapling.id values in column A: the entries that had values
# 152, 137 and 156 are moved from their original positions to the end of
# the list (rows 10, 36 and 37 respectively), shifting the rest of the
# rows up by one accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    222, 264, 297, 343, 442, 455, 500, 677, 689, 695,
    731, 753, 761, 1245, 1392, 1447, 1529, 25, 30, 70,
    86, 95, 99, 117, 119, 189, 449, 679, 682, 704,
    1473, 1475, 74, 142, 152, 137, 156
)

$startRow = 10
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
